# Update kd of NRP1 fitted data excel sheet
#
# - Refresh the G12:G16 raw kon readings on Sheet2 (new fit values), which
#   ripple into the STDEV-based SE formulas already on the sheet (H12, K12:K16, L12).
# - Switch the Kd formulas in Sheet1 J2:J4 from H/F (single replicate ratio)
#   to AVERAGE(I)/AVERAGE(G) pulled straight from Sheet2 (breaks the shared
#   formula group that used to link J3/J4).
# - Reposition the on-screen selection/scroll for both sheets the way the
#   author last left them before saving.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: updated kon replicate readings (column G, rows 12-16) ---
$ws2.Range("G12").Value = 87144780
$ws2.Range("G13").Value = 45296290
$ws2.Range("G14").Value = 24892240
$ws2.Range("G15").Value = 14307420
$ws2.Range("G16").Value = 7196948

# --- Sheet1: Kd formulas now average across the five Sheet2 replicates ---
$ws1.Range("J2").Formula = "=AVERAGE(Sheet2!I2:I6)/AVERAGE(Sheet2!G2:G6)"
$ws1.Range("J3").Formula = "=AVERAGE(Sheet2!I7:I11)/AVERAGE(Sheet2!G7:G11)"
$ws1.Range("J4").Formula = "=AVERAGE(Sheet2!I12:I16)/AVERAGE(Sheet2!G12:G16)"

# --- Restore the view state (active cell / scroll) left by the author ---
$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws2.Range("G13").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("I24").Select()
